$wb = $excel.ActiveWorkbook
$wsTransform = $wb.Worksheets.Item("transform")
$wsOutput = $wb.Worksheets.Item("output")

$wsTransform.Range("Q17").Value = "Corinna Blake, Marsha May, Romi Rain"
$wsOutput.Range("U17").Value = "Corinna Blake, Marsha May, Romi Rain"
$wsTransform.Range("Q22").Value = "Ash Hollywood, Juelz Ventura, Adessa Winters, Ariella Ferrera, London Keyes, Jessica Jaymes, Ava Addams, Keisha Grey"
$wsOutput.Range("U22").Value = "Ash Hollywood, Juelz Ventura, Adessa Winters, Ariella Ferrera, London Keyes, Jessica Jaymes, Ava Addams, Keisha Grey"
$wsTransform.Range("Q27").Value = "Ash Hollywood, Juelz Ventura, Adessa Winters, Ariella Ferrera, London Keyes, Jessica Jaymes, Ava Addams, Keisha Grey"
$wsOutput.Range("U27").Value = "Ash Hollywood, Juelz Ventura, Adessa Winters, Ariella Ferrera, London Keyes, Jessica Jaymes, Ava Addams, Keisha Grey"
$wsTransform.Range("Q32").Value = "Ash Hollywood, Juelz Ventura, Adessa Winters, Ariella Ferrera, London Keyes, Jessica Jaymes, Ava Addams, Keisha Grey"
$wsOutput.Range("U32").Value = "Ash Hollywood, Juelz Ventura, Adessa Winters, Ariella Ferrera, London Keyes, Jessica Jaymes, Ava Addams, Keisha Grey"
$wsTransform.Range("Q36").Value = "Remy LaCroix, Jada Stevens, Maddy Oreilly, Kennedy Leigh, Abby Cross, Mary Jane Mayhem, Rahyndee James"
$wsOutput.Range("U36").Value = "Remy LaCroix, Jada Stevens, Maddy Oreilly, Kennedy Leigh, Abby Cross, Mary Jane Mayhem, Rahyndee James"
$wsTransform.Range("Q40").Value = "Bianca Breeze, Charlotte Cross, Darla Crane, Gwen Stark, Amanda Lane, Veronica Rayne"
$wsOutput.Range("U40").Value = "Bianca Breeze, Charlotte Cross, Darla Crane, Gwen Stark, Amanda Lane, Veronica Rayne"
$wsTransform.Range("Q45").Value = "Abigail Mac, Ashley Fires, JoJo Kiss, Alessa Savage, Aylin Diamond"
$wsOutput.Range("U45").Value = "Abigail Mac, Ashley Fires, JoJo Kiss, Alessa Savage, Aylin Diamond"
$wsTransform.Range("Q50").Value = "Abigail Mac, Ashley Fires, JoJo Kiss, Alessa Savage, Aylin Diamond"
$wsOutput.Range("U50").Value = "Abigail Mac, Ashley Fires, JoJo Kiss, Alessa Savage, Aylin Diamond"
$wsTransform.Range("Q54").Value = "Alexis Fawx, Cassidy Klein, Simone Sonay, Carmel Anderson, Sensual Jane "
$wsOutput.Range("U54").Value = "Alexis Fawx, Cassidy Klein, Simone Sonay, Carmel Anderson, Sensual Jane "
$wsTransform.Range("Q59").Value = "Allie Haze, Peta Jensen, Romi Rain, Allie Haze, Romi Rain, Romi Rain"
$wsOutput.Range("U59").Value = "Allie Haze, Peta Jensen, Romi Rain, Allie Haze, Romi Rain, Romi Rain"
$wsTransform.Range("Q64").Value = "Allie Haze, Peta Jensen, Romi Rain, Allie Haze, Romi Rain, Romi Rain"
$wsOutput.Range("U64").Value = "Allie Haze, Peta Jensen, Romi Rain, Allie Haze, Romi Rain, Romi Rain"
$wsTransform.Range("Q69").Value = "Allie Haze, Peta Jensen, Romi Rain, Allie Haze, Romi Rain, Romi Rain"
$wsOutput.Range("U69").Value = "Allie Haze, Peta Jensen, Romi Rain, Allie Haze, Romi Rain, Romi Rain"
$wsTransform.Range("Q73").Value = "Isis Love, Peta Jensen, Kissa Sins, Pristine Edge"
$wsOutput.Range("U73").Value = "Isis Love, Peta Jensen, Kissa Sins, Pristine Edge"
$wsTransform.Range("Q77").Value = "Marica Hase, Sharon Lee, Jayden Lee"
$wsOutput.Range("U77").Value = "Marica Hase, Sharon Lee, Jayden Lee"
$wsTransform.Range("Q81").Value = "Marica Hase, Sharon Lee, Jayden Lee"
$wsOutput.Range("U81").Value = "Marica Hase, Sharon Lee, Jayden Lee"
$wsTransform.Range("Q85").Value = "Kimmy and Marsha Have a Threesome, Cory Chase, Billi Ann, Elsa Jean, Valentinna Nappi"
$wsOutput.Range("U85").Value = "Kimmy and Marsha Have a Threesome, Cory Chase, Billi Ann, Elsa Jean, Valentinna Nappi"
$wsTransform.Range("Q89").Value = "Gia Paige, Kimber Delice, Felicia Kiss"
$wsOutput.Range("U89").Value = "Gia Paige, Kimber Delice, Felicia Kiss"
$wsTransform.Range("Q93").Value = "Gia Paige, Kimber Delice, Felicia Kiss"
$wsOutput.Range("U93").Value = "Gia Paige, Kimber Delice, Felicia Kiss"
$wsTransform.Range("Q98").Value = "Ryta, Antonia Sainz, Kimber Delice, Daphne Klyde"
$wsOutput.Range("U98").Value = "Ryta, Antonia Sainz, Kimber Delice, Daphne Klyde"
$wsTransform.Range("Q103").Value = "Ryta, Antonia Sainz, Kimber Delice, Daphne Klyde"
$wsOutput.Range("U103").Value = "Ryta, Antonia Sainz, Kimber Delice, Daphne Klyde"
$wsTransform.Range("Q108").Value = "Cassidy Banks, Yurizan Beltran, Stacy Jay, Cristi Ann"
$wsOutput.Range("U108").Value = "Cassidy Banks, Yurizan Beltran, Stacy Jay, Cristi Ann"
$wsTransform.Range("Q113").Value = "Cassidy Banks, Yurizan Beltran, Stacy Jay, Cristi Ann"
$wsOutput.Range("U113").Value = "Cassidy Banks, Yurizan Beltran, Stacy Jay, Cristi Ann"
$wsTransform.Range("Q118").Value = "CeCe Capella, Lilly Sapphire, Lucy Doll, Brittany Shae"
$wsOutput.Range("U118").Value = "CeCe Capella, Lilly Sapphire, Lucy Doll, Brittany Shae"
$wsTransform.Range("Q123").Value = "CeCe Capella, Lilly Sapphire, Lucy Doll, Brittany Shae"
$wsOutput.Range("U123").Value = "CeCe Capella, Lilly Sapphire, Lucy Doll, Brittany Shae"
$wsTransform.Range("Q127").Value = "Ally, Ava Alba, Melissa Moore"
$wsOutput.Range("U127").Value = "Ally, Ava Alba, Melissa Moore"
$wsTransform.Range("Q131").Value = "Ally, Ava Alba, Melissa Moore"
$wsOutput.Range("U131").Value = "Ally, Ava Alba, Melissa Moore"
$wsTransform.Range("Q136").Value = "Abigail Mac, Aspen Rae, Cayla Lyons, Naomi Nevena, Maci Winslett, Staci Carr, Lucy Li, Vinna Reed"
$wsOutput.Range("U136").Value = "Abigail Mac, Aspen Rae, Cayla Lyons, Naomi Nevena, Maci Winslett, Staci Carr, Lucy Li, Vinna Reed"
$wsTransform.Range("Q141").Value = "Abigail Mac, Aspen Rae, Cayla Lyons, Naomi Nevena, Maci Winslett, Staci Carr, Lucy Li, Vinna Reed"
$wsOutput.Range("U141").Value = "Abigail Mac, Aspen Rae, Cayla Lyons, Naomi Nevena, Maci Winslett, Staci Carr, Lucy Li, Vinna Reed"
$wsTransform.Range("Q145").Value = "Uma Jolie, Naomi Woods, Monica Belluci"
$wsOutput.Range("U145").Value = "Uma Jolie, Naomi Woods, Monica Belluci"
$wsTransform.Range("Q154").Value = "Aaliyah Love, Taylor Sands, Ava Dalush, Bianca Breeze"
$wsOutput.Range("U154").Value = "Aaliyah Love, Taylor Sands, Ava Dalush, Bianca Breeze"
$wsTransform.Range("Q158").Value = "Kissa Sins, Dani Daniels, Jessa Rhodes, Kayla Carrera, Kendra James"
$wsOutput.Range("U158").Value = "Kissa Sins, Dani Daniels, Jessa Rhodes, Kayla Carrera, Kendra James"
$wsTransform.Range("Q167").Value = "Jean Michaels, August Ames, Dillion Harper, Tia Cyrus"
$wsOutput.Range("U167").Value = "Jean Michaels, August Ames, Dillion Harper, Tia Cyrus"
$wsTransform.Range("Q171").Value = "Sara Jay, Makayla Cox, Richelle Ryan"
$wsOutput.Range("U171").Value = "Sara Jay, Makayla Cox, Richelle Ryan"
$wsTransform.Range("Q175").Value = "Sara Jay, Makayla Cox, Richelle Ryan"
$wsOutput.Range("U175").Value = "Sara Jay, Makayla Cox, Richelle Ryan"
$wsTransform.Range("Q179").Value = "Lea Lexis, Katy Rose, Leanna Sweet, Nekane, Billie Star, Abrill Gerald"
$wsOutput.Range("U179").Value = "Lea Lexis, Katy Rose, Leanna Sweet, Nekane, Billie Star, Abrill Gerald"
$wsTransform.Range("Q184").Value = "Nikki Benz, Jayden Jaymes, Delta White, Nikki Benz, Nikki Benz, Nikki Sexx"
$wsOutput.Range("U184").Value = "Nikki Benz, Jayden Jaymes, Delta White, Nikki Benz, Nikki Benz, Nikki Sexx"
$wsTransform.Range("Q189").Value = "Rachel Starr, Jenni Lee, Asa Akira, Mia Malkova"
$wsOutput.Range("U189").Value = "Rachel Starr, Jenni Lee, Asa Akira, Mia Malkova"
$wsTransform.Range("Q200").Value = "Keisha Grey, Mia Malkova, Lizz Taylor, Veronica Rodriguez, Rachel Roxxx, Chirsty Mack, Brandi Love, Riley Reid, Nikki Benz, Amy Anderssen, Courtney Taylor, Summer Brielle, Raven Bay, Rikki Six, Madison Ivy, Monique Alexander, Courtney Cummz, Tasha Reign, Brooklyn Chase, Kagney Linn Karter"
$wsOutput.Range("U200").Value = "Keisha Grey, Mia Malkova, Lizz Taylor, Veronica Rodriguez, Rachel Roxxx, Chirsty Mack, Brandi Love, Riley Reid, Nikki Benz, Amy Anderssen, Courtney Taylor, Summer Brielle, Raven Bay, Rikki Six, Madison Ivy, Monique Alexander, Courtney Cummz, Tasha Reign, Brooklyn Chase, Kagney Linn Karter"
$wsTransform.Range("Q211").Value = "Darling Danika, Katrina Jade, Bonnie Rotten, Mia Malkova, Audrey Bitoni, Richelle Ryan, Asa Akira, Christy Mack, Summer Brielle, Madison Scott, Kagney Linn Karter, Cherie Deville, Siri"
$wsOutput.Range("U211").Value = "Darling Danika, Katrina Jade, Bonnie Rotten, Mia Malkova, Audrey Bitoni, Richelle Ryan, Asa Akira, Christy Mack, Summer Brielle, Madison Scott, Kagney Linn Karter, Cherie Deville, Siri"
$wsTransform.Range("Q222").Value = "Darla Crane, Devon, Veronica Avluv, Priya Anjali Rai, Lisa Ann, Tiffany Mynx, Nikita Von James, Raylene, Deauxma, Houston"
$wsOutput.Range("U222").Value = "Darla Crane, Devon, Veronica Avluv, Priya Anjali Rai, Lisa Ann, Tiffany Mynx, Nikita Von James, Raylene, Deauxma, Houston"
$wsTransform.Range("Q233").Value = "Marica Hase, Akira Lane, Katsuni, Alina Li, Jayden Lee, Harumi Asano, Kianna Dior, Asa Akira, London Keyes, Kalina Ryu, Morgan Lee, Kaylani Lei"
$wsOutput.Range("U233").Value = "Marica Hase, Akira Lane, Katsuni, Alina Li, Jayden Lee, Harumi Asano, Kianna Dior, Asa Akira, London Keyes, Kalina Ryu, Morgan Lee, Kaylani Lei"
$wsTransform.Range("Q238").Value = "Alex Mecum, Dennis West, Dennis West, Luke Adams, Will Braun, Dennis West, Landon Mycles, Dennis West, Scott Riley"
$wsOutput.Range("U238").Value = "Alex Mecum, Dennis West, Dennis West, Luke Adams, Will Braun, Dennis West, Landon Mycles, Dennis West, Scott Riley"
$wsTransform.Range("Q243").Value = "Alex Mecum, Dennis West, Dennis West, Luke Adams, Will Braun, Dennis West, Landon Mycles, Dennis West, Scott Riley"
$wsOutput.Range("U243").Value = "Alex Mecum, Dennis West, Dennis West, Luke Adams, Will Braun, Dennis West, Landon Mycles, Dennis West, Scott Riley"
$wsTransform.Range("Q248").Value = "Alex Mecum, Dennis West, Dennis West, Luke Adams, Will Braun, Dennis West, Landon Mycles, Dennis West, Scott Riley"
$wsOutput.Range("U248").Value = "Alex Mecum, Dennis West, Dennis West, Luke Adams, Will Braun, Dennis West, Landon Mycles, Dennis West, Scott Riley"
$wsTransform.Range("Q253").Value = "Phenix Saint, Tommy Regan, Dennis West, Jake Bass, Phenix Saint, Scott Riley, Jimmy Fanz, Phenix Saint"
$wsOutput.Range("U253").Value = "Phenix Saint, Tommy Regan, Dennis West, Jake Bass, Phenix Saint, Scott Riley, Jimmy Fanz, Phenix Saint"
$wsTransform.Range("Q258").Value = "Phenix Saint, Tommy Regan, Dennis West, Jake Bass, Phenix Saint, Scott Riley, Jimmy Fanz, Phenix Saint"
$wsOutput.Range("U258").Value = "Phenix Saint, Tommy Regan, Dennis West, Jake Bass, Phenix Saint, Scott Riley, Jimmy Fanz, Phenix Saint"
$wsTransform.Range("Q263").Value = "Phenix Saint, Tommy Regan, Dennis West, Jake Bass, Phenix Saint, Scott Riley, Jimmy Fanz, Phenix Saint"
$wsOutput.Range("U263").Value = "Phenix Saint, Tommy Regan, Dennis West, Jake Bass, Phenix Saint, Scott Riley, Jimmy Fanz, Phenix Saint"
$wsTransform.Range("Q268").Value = "Blue Moores, Paddy O'Brian, Dimitri Kane, Tony Paradise, Landon Mycles, Sebastian Young, Adam Bryant, Josh Peters"
$wsOutput.Range("U268").Value = "Blue Moores, Paddy O'Brian, Dimitri Kane, Tony Paradise, Landon Mycles, Sebastian Young, Adam Bryant, Josh Peters"
$wsTransform.Range("Q273").Value = "Blue Moores, Paddy O'Brian, Dimitri Kane, Tony Paradise, Landon Mycles, Sebastian Young, Adam Bryant, Josh Peters"
$wsOutput.Range("U273").Value = "Blue Moores, Paddy O'Brian, Dimitri Kane, Tony Paradise, Landon Mycles, Sebastian Young, Adam Bryant, Josh Peters"
$wsTransform.Range("Q278").Value = "Peter Fields, Will Braun, Luke Adams, Will Braun, Damien Michaels, Jordan Boss, Josh Peters, Will Braun"
$wsOutput.Range("U278").Value = "Peter Fields, Will Braun, Luke Adams, Will Braun, Damien Michaels, Jordan Boss, Josh Peters, Will Braun"
$wsTransform.Range("Q283").Value = "Peter Fields, Will Braun, Luke Adams, Will Braun, Damien Michaels, Jordan Boss, Josh Peters, Will Braun"
$wsOutput.Range("U283").Value = "Peter Fields, Will Braun, Luke Adams, Will Braun, Damien Michaels, Jordan Boss, Josh Peters, Will Braun"
$wsTransform.Range("Q288").Value = "Peter Fields, Will Braun, Luke Adams, Will Braun, Damien Michaels, Jordan Boss, Josh Peters, Will Braun"
$wsOutput.Range("U288").Value = "Peter Fields, Will Braun, Luke Adams, Will Braun, Damien Michaels, Jordan Boss, Josh Peters, Will Braun"
$wsTransform.Range("Q293").Value = "Jack Hunter, Will Braun, Colby Keller, Will Braun, Colby Keller, Roman Todd, Addison Graham, Brandon Moore, Colby Keller, Roman Todd, Will Braun"
$wsOutput.Range("U293").Value = "Jack Hunter, Will Braun, Colby Keller, Will Braun, Colby Keller, Roman Todd, Addison Graham, Brandon Moore, Colby Keller, Roman Todd, Will Braun"
$wsTransform.Range("Q298").Value = "Jack Hunter, Will Braun, Colby Keller, Will Braun, Colby Keller, Roman Todd, Addison Graham, Brandon Moore, Colby Keller, Roman Todd, Will Braun"
$wsOutput.Range("U298").Value = "Jack Hunter, Will Braun, Colby Keller, Will Braun, Colby Keller, Roman Todd, Addison Graham, Brandon Moore, Colby Keller, Roman Todd, Will Braun"
$wsTransform.Range("Q303").Value = "Diego Sans, Jordan Boss, Diego Sans, Ricky Decker, Bennett Anthony, Roman Todd, Diego Sans, Scott Riley"
$wsOutput.Range("U303").Value = "Diego Sans, Jordan Boss, Diego Sans, Ricky Decker, Bennett Anthony, Roman Todd, Diego Sans, Scott Riley"
$wsTransform.Range("Q308").Value = "Diego Sans, Jordan Boss, Diego Sans, Ricky Decker, Bennett Anthony, Roman Todd, Diego Sans, Scott Riley"
$wsOutput.Range("U308").Value = "Diego Sans, Jordan Boss, Diego Sans, Ricky Decker, Bennett Anthony, Roman Todd, Diego Sans, Scott Riley"
$wsTransform.Range("Q313").Value = "Diego Sans, Jordan Boss, Diego Sans, Ricky Decker, Bennett Anthony, Roman Todd, Diego Sans, Scott Riley"
$wsOutput.Range("U313").Value = "Diego Sans, Jordan Boss, Diego Sans, Ricky Decker, Bennett Anthony, Roman Todd, Diego Sans, Scott Riley"
$wsTransform.Range("Q318").Value = "Damien Crosse, Dario Beck, Colby Jansen, Dario Beck, Dario Beck, Massimo Piano, JJ Knight, Tommy Regan"
$wsOutput.Range("U318").Value = "Damien Crosse, Dario Beck, Colby Jansen, Dario Beck, Dario Beck, Massimo Piano, JJ Knight, Tommy Regan"
$wsTransform.Range("Q323").Value = "Damien Crosse, Dario Beck, Colby Jansen, Dario Beck, Dario Beck, Massimo Piano, JJ Knight, Tommy Regan"
$wsOutput.Range("U323").Value = "Damien Crosse, Dario Beck, Colby Jansen, Dario Beck, Dario Beck, Massimo Piano, JJ Knight, Tommy Regan"
$wsTransform.Range("Q328").Value = "Alex Mecum, Aspen, Klein Kerr, Massimo Piano, Bruno Bernal, Jake Bass, Alex Mecum, Colby Keller"
$wsOutput.Range("U328").Value = "Alex Mecum, Aspen, Klein Kerr, Massimo Piano, Bruno Bernal, Jake Bass, Alex Mecum, Colby Keller"
$wsTransform.Range("Q333").Value = "Alex Mecum, Aspen, Klein Kerr, Massimo Piano, Bruno Bernal, Jake Bass, Alex Mecum, Colby Keller"
$wsOutput.Range("U333").Value = "Alex Mecum, Aspen, Klein Kerr, Massimo Piano, Bruno Bernal, Jake Bass, Alex Mecum, Colby Keller"
$wsTransform.Range("Q338").Value = "Alex Mecum, Aspen, Klein Kerr, Massimo Piano, Bruno Bernal, Jake Bass, Alex Mecum, Colby Keller"
$wsOutput.Range("U338").Value = "Alex Mecum, Aspen, Klein Kerr, Massimo Piano, Bruno Bernal, Jake Bass, Alex Mecum, Colby Keller"
$wsTransform.Range("Q349").Value = "Adam Killian, John Magnum, Johnny Rapid, Sebastian Young, Christian Wilde, Cliff Jensen, Gavin Waters, Phillip Aubrey, Adam Killian, Trenton Ducati, Colby Keller, Jake Steel, Dereck Fox, Tyler St.James, Boston Miles, Johnny Rapid, Travis Irons , Colby Jansen, Duncan Black, Parker London, Phenix Saint, Sebastian Keys, Tony Paradise, Rafael Alencar, Johnny Rapid, Spencer Fox"
$wsOutput.Range("U349").Value = "Adam Killian, John Magnum, Johnny Rapid, Sebastian Young, Christian Wilde, Cliff Jensen, Gavin Waters, Phillip Aubrey, Adam Killian, Trenton Ducati, Colby Keller, Jake Steel, Dereck Fox, Tyler St.James, Boston Miles, Johnny Rapid, Travis Irons , Colby Jansen, Duncan Black, Parker London, Phenix Saint, Sebastian Keys, Tony Paradise, Rafael Alencar, Johnny Rapid, Spencer Fox"
$wsTransform.Range("Q360").Value = "Adam Killian, John Magnum, Johnny Rapid, Sebastian Young, Christian Wilde, Cliff Jensen, Gavin Waters, Phillip Aubrey, Adam Killian, Trenton Ducati, Colby Keller, Jake Steel, Dereck Fox, Tyler St.James, Boston Miles, Johnny Rapid, Travis Irons , Colby Jansen, Duncan Black, Parker London, Phenix Saint, Sebastian Keys, Tony Paradise, Rafael Alencar, Johnny Rapid, Spencer Fox"
$wsOutput.Range("U360").Value = "Adam Killian, John Magnum, Johnny Rapid, Sebastian Young, Christian Wilde, Cliff Jensen, Gavin Waters, Phillip Aubrey, Adam Killian, Trenton Ducati, Colby Keller, Jake Steel, Dereck Fox, Tyler St.James, Boston Miles, Johnny Rapid, Travis Irons , Colby Jansen, Duncan Black, Parker London, Phenix Saint, Sebastian Keys, Tony Paradise, Rafael Alencar, Johnny Rapid, Spencer Fox"

# Update media_sourced_from timestamp string in output sheet column C (rows 2-482)
$wsOutput.Range("C2:C482").Value = "02/18/2016 10:29:32"

